$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: reorder "Recorded By" list for ANATOMY session 1
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System"

# G9: reorder "Recorded By" list for HISTOLOGY session 1
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# L10: Average Attendance % for HISTOLOGY group stats (stored as plain text,
# not a numeric percentage) - use a leading apostrophe so it is entered as
# text instead of being auto-converted to a percentage number, then restore
# the original cell formatting (General, no quote-prefix) by copying the
# format from the untouched neighboring K10 cell which shares the same style.
$ws.Range("L10").Value = "'26.5%"
$ws.Range("K10").Copy() | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null

# H14: Students attendance for PARASITOLOGY session 1
$ws.Range("H14").Value = "93/251"

# S15: Avg Attendance % for PARASITOLOGY group stats (same text treatment as L10)
$ws.Range("S15").Value = "'26.5%"
$ws.Range("K10").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null

# G28: reorder "Recorded By" list for PHYSIOLOGY session 1
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

$excel.CutCopyMode = $false
